# Test_PEA_verifierAccesBadge.docx edit script
# - Merge "verifierAcces"/"Badge" spell-checked runs (handled automatically by the
#   engine's XML round-trip/normalisation when we re-insert the document XML).
# - Remove now-stray <w:proofErr/> markers around "pea"/"acces" (also handled by
#   the round-trip normalisation).
# - Resize the results table's columns.
# - Clear the "Statut attendu" header cell and all of the "Erreur"/"Succès"
#   status cells in the results table, and drop their <w:hideMark/> hint.

$d = $word.ActiveDocument

$xml = $d.Content.WordOpenXML

# --- Resize the grid columns of the results table -------------------------------
$xml = $xml.Replace('<w:gridCol w:w="587"/>',  '<w:gridCol w:w="621"/>')
$xml = $xml.Replace('<w:gridCol w:w="2093"/>', '<w:gridCol w:w="2328"/>')
$xml = $xml.Replace('<w:gridCol w:w="2939"/>', '<w:gridCol w:w="3406"/>')
$xml = $xml.Replace('<w:gridCol w:w="2352"/>', '<w:gridCol w:w="2636"/>')
$xml = $xml.Replace('<w:gridCol w:w="1101"/>', '<w:gridCol w:w="81"/>')

# --- Clear the "Statut attendu" header cell (keeps the (now empty) paragraph) ---
$xml = $xml.Replace(
    '<w:tcPr><w:tcW w:w="0" w:type="auto"/><w:vAlign w:val="center"/><w:hideMark/></w:tcPr><w:p w14:paraId="7CED20D8" w14:textId="77777777" w:rsidR="005E0781" w:rsidRPr="005E0781" w:rsidRDefault="005E0781" w:rsidP="005E0781"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="005E0781"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Statut attendu</w:t></w:r></w:p>',
    '<w:tcPr><w:tcW w:w="0" w:type="auto"/><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="7CED20D8" w14:textId="77777777" w:rsidR="005E0781" w:rsidRPr="005E0781" w:rsidRDefault="005E0781" w:rsidP="005E0781"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p>'
)

# --- Clear the "Erreur"/"Succès" status cells -----------------------------------
$statusParaIds = @(
    "5BD3D637","0B3BC608","0F2008C5","012F9755","31444445",
    "1767328B","58D38ED4","15B4C4A3","16F8C3A7","48AD7D69","26232E67"
)

foreach ($pid in $statusParaIds) {
    $pattern = '<w:tcPr><w:tcW w:w="0" w:type="auto"/><w:vAlign w:val="center"/><w:hideMark/></w:tcPr><w:p w14:paraId="' + $pid + '" w14:textId="77777777" w:rsidR="005E0781" w:rsidRPr="005E0781" w:rsidRDefault="005E0781" w:rsidP="005E0781"><w:r w:rsidRPr="005E0781"><w:t>Erreur</w:t></w:r></w:p>'
    $replacement = '<w:tcPr><w:tcW w:w="0" w:type="auto"/><w:vAlign w:val="center"/></w:tcPr><w:p/>'
    $xml = $xml.Replace($pattern, $replacement)

    $pattern2 = '<w:tcPr><w:tcW w:w="0" w:type="auto"/><w:vAlign w:val="center"/><w:hideMark/></w:tcPr><w:p w14:paraId="' + $pid + '" w14:textId="77777777" w:rsidR="005E0781" w:rsidRPr="005E0781" w:rsidRDefault="005E0781" w:rsidP="005E0781"><w:r w:rsidRPr="005E0781"><w:t>Succès</w:t></w:r></w:p>'
    $xml = $xml.Replace($pattern2, $replacement)
}

$d.Content.InsertXML($xml)
